# Commit: "Change TestObject folder structure and page naming convention
# to separate from functional tests"
#
# The sheets that back the Bento static-data workbook are renamed with a
# "V_" prefix (and the underscore in "Home_page" is dropped), the
# previously-selected tab ("Home_page") loses focus, and "AboutBentoPage"
# becomes the new active/selected tab. Along with that tab switch, the
# remembered selections on the two affected sheets move to new cells.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets (folder / naming convention change) -----------------
$wb.Worksheets.Item("Home_page").Name          = "V_HomePage"
$wb.Worksheets.Item("AboutBentoPage").Name     = "V_AboutBentoPage"
$wb.Worksheets.Item("AboutResourcesPage").Name = "V_AboutResourcesPage"

$wsHome  = $wb.Worksheets.Item("V_HomePage")
$wsBento = $wb.Worksheets.Item("V_AboutBentoPage")

# --- 2. V_HomePage: no longer the selected tab, selection moves to A40 ----
$wsHome.Activate() | Out-Null
$wsHome.Range("A40").Select() | Out-Null

# --- 3. V_AboutBentoPage: becomes the selected/active tab, selection -> A6 -
$wsBento.Activate() | Out-Null
$wsBento.Range("A6").Select() | Out-Null
